$wb = $excel.ActiveWorkbook

# Sheet 1: "VENTAS POR GRUPO" -> H22 changes from -67.65000000000001 to -255.75
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("H22").Value = -255.75

# Sheet 2: "VENTA MENSUAL" -> F22 changes from 7229.68 to 7041.58
#                               F23 changes from 24077.51 to 23889.41
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F22").Value = 7041.58
$ws2.Range("F23").Value = 23889.41

# Sheet 3: "CUMPLIMIENTO MENSUAL"
#   Row 6 (INODOROS): D6 1735.95 -> 1547.85 ; E6 1171.63368146026 -> 1359.73368146026 ; F6 0.5970421457064181 -> 0.5323492527040982
#   Row 15 (TOTAL):   D15 24077.51 -> 23889.41 ; E15 31347.23316613378 -> 31535.33316613378 ; F15 0.4344180707852535 -> 0.431024279686643
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D6").Value = 1547.85
$ws3.Range("E6").Value = 1359.73368146026
$ws3.Range("F6").Value = 0.5323492527040982

$ws3.Range("D15").Value = 23889.41
$ws3.Range("E15").Value = 31535.33316613378
$ws3.Range("F15").Value = 0.431024279686643
